$d = $word.ActiveDocument

# Locate the existing sentence that needs the " (Changed main)" suffix
# split across separate runs, matching the target OOXML exactly.
$target = "This is a Microsoft word document."
$findRange = $d.Content
$found = $findRange.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found) {
    throw "Could not find target sentence '$target' in the document."
}

# $findRange now spans exactly the matched text (no paragraph mark included).
$editRange = $d.Range($findRange.Start, $findRange.End)

# Rebuild the paragraph's run content as four separate runs:
#   1) the original sentence (unchanged)
#   2) " ("
#   3) "Changed main"
#   4) ")"
# InsertXML replaces the contents of the addressed range, so by targeting
# exactly the original run's text we keep the rest of the paragraph (and
# document) untouched while forcing genuinely separate <w:r> elements.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>This is a Microsoft word document.</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> (</w:t>
            </w:r>
            <w:r>
              <w:t>Changed main</w:t>
            </w:r>
            <w:r>
              <w:t>)</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$editRange.InsertXML($xml)
